$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.086.80'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.480.48'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.93'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.69'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.513'
$ws.Range('E8').Value = '  -1.21%  '
$ws.Range('D9').Value = '2.480.66'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.165'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.93'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('D14').Value = '2.913.00'
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.30'
$ws.Range('E15').Value = '  -2.45%  '
$ws.Range('D16').Value = '66.999.75'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').Value = '2.503.17'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.94'
$ws.Range('E19').Value = '  -5.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.39'
$ws.Range('E20').Value = '  -6.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.98'
$ws.Range('E21').Value = '  -3.64%  '
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.55'
$ws.Range('E24').Value = '  -4.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.22'
$ws.Range('E25').Value = '  -5.45%  '
$ws.Range('E26').Value = '  -2.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.15'
$ws.Range('E27').Value = '  -3.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -9.13%  '
$ws.Range('D29').Value = '2.608.26'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('D30').Value = '0.0₃0901'
$ws.Range('E30').Value = '  -3.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '507.84'
$ws.Range('E31').Value = '  -2.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.67'
$ws.Range('E32').Value = '  -6.04%  '
$ws.Range('E33').Value = '  -4.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.23'
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.117'
$ws.Range('E36').Value = '  -7.46%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.72'
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.68'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.21'
$ws.Range('E39').Value = '  -4.41%  '
$ws.Range('E40').Value = '  -6.40%  '
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('E42').Value = '  -3.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.327'
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('E44').Value = '  -4.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.37'
$ws.Range('E45').Value = '  -4.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.72'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.10'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.514'
$ws.Range('E48').Value = '  -4.43%  '
$ws.Range('E49').Value = '  -5.26%  '
$ws.Range('E50').Value = '  -6.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0729'
$ws.Range('E51').Value = '  -1.11%  '
